# Sources_overview.xlsx clean-up
# - Replace the raw source-URL / ad-hoc note columns with a cleaner
#   "Source" (authority name) + "Procured via" (website/email) layout.
# - Re-point the remaining hyperlinks at the authorities' homepages.
# - Resize columns B/C to fit the new, shorter text.
# - Move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Drop every pre-existing hyperlink (Hessen / Niedersachsen lose theirs
#    entirely; Schleswig-Holstein / Thueringen get new targets below).
# ---------------------------------------------------------------------
$ws.Range("A1").Hyperlinks.Delete()

# ---------------------------------------------------------------------
# 2) Rewrite the table contents (columns A-C, rows 1-17).
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Bundesland"
$ws.Range("B1").Value = "Source"
$ws.Range("C1").Value = "Procured via"

$ws.Range("A2").Value = "Baden-Wuerttemberg"
$ws.Range("B2").Value = "Statistisches Landesamt Baden-Württemberg"
$ws.Range("C2").Value = "email"

$ws.Range("A3").Value = "Bayern"
$ws.Range("B3").Value = "Bayerisches Landesamt für Statistik"
$ws.Range("C3").Value = "website"

$ws.Range("A4").Value = "Brandenburg"
$ws.Range("B4").Value = "Amt für Statistik Berlin-Brandenburg"
$ws.Range("C4").Value = "website"

$ws.Range("A5").Value = "Bremen"
$ws.Range("B5").Value = "Statistisches Landesamt Bremen"
$ws.Range("C5").Value = "website"

$ws.Range("A6").Value = "BW"
$ws.Range("B6").Value = "Statistisches Landesamt Baden-Württemberg"
$ws.Range("C6").Value = "email"

$ws.Range("A7").Value = "Hamburg"
$ws.Range("B7").Value = "Statistik Nord"
$ws.Range("C7").Value = "website"

$ws.Range("A8").Value = "Hessen"
$ws.Range("B8").Value = "Hessisches Statistisches Landesamt"
$ws.Range("C8").Value = "website"

$ws.Range("A9").Value = "Mecklenburg Vorpommern"
$ws.Range("B9").Value = "Mecklenburg-Vorpommern Landesamt für innere Verwaltung Statistisches Amt"
$ws.Range("C9").Value = "website"

$ws.Range("A10").Value = "Niedersachsen"
$ws.Range("B10").Value = "Landesamt für Statistik Niedersachsen"
$ws.Range("C10").Value = "website after 2006; email before 2006"

$ws.Range("A11").Value = "NRW"
$ws.Range("B11").Value = "Statistisches Landesamt Nordrhein-Westfalen"
$ws.Range("C11").Value = "email"

$ws.Range("A12").Value = "RLP"
$ws.Range("B12").Value = "Statistisches Landesamt Rheinland-Pfalz"
$ws.Range("C12").Value = "email"

$ws.Range("A13").Value = "Saarland"
$ws.Range("B13").Value = "Statistisches Landesamt des Saarlandes"
$ws.Range("C13").Value = "email"

$ws.Range("A14").Value = "Sachsen"
$ws.Range("B14").Value = "Statistisches Landesamt des Freistaates Sachsen"
$ws.Range("C14").Value = "website"

$ws.Range("A15").Value = "Sachsen-Anhalt"
$ws.Range("B15").Value = "Statistisches Landesamt Sachsen-Anhalt"
$ws.Range("C15").Value = "website"

$ws.Range("A16").Value = "Schleswig-Holstein"
$ws.Range("B16").Value = "Statistisches Amt für Hamburg und Schleswig-Holstein"
$ws.Range("C16").Value = "website except for 2013; email for 2013"

$ws.Range("A17").Value = "Thueringen"
$ws.Range("B17").Value = "Thüringer Landesamt für Statistik"
$ws.Range("C17").Value = "website"

# ---------------------------------------------------------------------
# 3) Cells that used to carry the "Hyperlink" cell style but no longer
#    hold a link (Hessen, Niedersachsen) go back to the Normal style.
# ---------------------------------------------------------------------
$ws.Range("B8").Style = "Normal"
$ws.Range("B10").Style = "Normal"

# ---------------------------------------------------------------------
# 4) Add the (new / updated) hyperlinks. Hyperlinks.Add stamps the cell
#    text with the TextToDisplay argument, so re-apply the authority
#    name afterwards and strip the blue/underline "Hyperlink" style
#    Excel auto-applies, keeping the sheet's plain look.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("B2"), "https://www.statistik-bw.de/", "", "", "https://www.statistik-bw.de/")
$ws.Range("B2").Value = "Statistisches Landesamt Baden-Württemberg"
$ws.Range("B2").Style = "Normal"

$ws.Hyperlinks.Add($ws.Range("B11"), "https://www.it.nrw/statistik", "", "", "https://www.it.nrw/statistik")
$ws.Range("B11").Value = "Statistisches Landesamt Nordrhein-Westfalen"
$ws.Range("B11").Style = "Normal"

$ws.Hyperlinks.Add($ws.Range("B12"), "https://www.statistik.rlp.de/", "", "", "https://www.statistik.rlp.de/")
$ws.Range("B12").Value = "Statistisches Landesamt Rheinland-Pfalz"
$ws.Range("B12").Style = "Normal"

$ws.Hyperlinks.Add($ws.Range("B13"), "https://www.saarland.de/stat/DE/home", "", "", "https://www.saarland.de/stat/DE/home")
$ws.Range("B13").Value = "Statistisches Landesamt des Saarlandes"
$ws.Range("B13").Style = "Normal"

$ws.Hyperlinks.Add($ws.Range("B15"), "https://statistik.sachsen-anhalt.de/", "", "", "https://statistik.sachsen-anhalt.de/")
$ws.Range("B15").Value = "Statistisches Landesamt Sachsen-Anhalt"
$ws.Range("B15").Style = "Normal"

$ws.Hyperlinks.Add($ws.Range("B16"), "https://www.statistik-nord.de/", "", "", "https://www.statistik-nord.de/")
$ws.Range("B16").Value = "Statistisches Amt für Hamburg und Schleswig-Holstein"
$ws.Range("B16").Style = "Normal"

$ws.Hyperlinks.Add($ws.Range("B17"), "https://statistik.thueringen.de/", "", "", "https://statistik.thueringen.de/")
$ws.Range("B17").Value = "Thüringer Landesamt für Statistik"
$ws.Range("B17").Style = "Normal"

# ---------------------------------------------------------------------
# 5) Column widths: B shrinks to fit the (shorter) authority names, the
#    new column C is sized for the "Procured via" text.
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 65.7109375
$ws.Columns.Item(3).ColumnWidth = 33.75

# ---------------------------------------------------------------------
# 6) Move the selection like the author's last click.
# ---------------------------------------------------------------------
$ws.Range("C18").Select()
